$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 403, pushing the existing
# rows 403:418 down to 405:420 (this also grows the used range /
# dimension to A1:R420 automatically).
$ws.Rows.Item(403).Insert()
$ws.Rows.Item(403).Insert()

# Populate the two newly-inserted rows (new weekly observations).
$ws.Range("A403").Value = 7
$ws.Range("B403").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C403").Value = "Ñuble"
$ws.Range("D403").Value = 44509
$ws.Range("E403").Value = 16
$ws.Range("F403").Value = 100112033
$ws.Range("G403").Value = "Lechuga"
$ws.Range("H403").Value = "Conconina(o)"
$ws.Range("I403").Value = "Primera"
$ws.Range("J403").Value = 240
$ws.Range("K403").Value = 5000
$ws.Range("L403").Value = 5500
$ws.Range("M403").Value = 5250
$ws.Range("N403").Value = "`$/caja 10 unidades"
$ws.Range("O403").Value = "Región del Maule"
$ws.Range("P403").Value = 525
$ws.Range("Q403").Value = 10
$ws.Range("R403").Value = "Hortaliza"

$ws.Range("A404").Value = 7
$ws.Range("B404").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C404").Value = "Ñuble"
$ws.Range("D404").Value = 44509
$ws.Range("E404").Value = 16
$ws.Range("F404").Value = 100112033
$ws.Range("G404").Value = "Lechuga"
$ws.Range("H404").Value = "Escarola"
$ws.Range("I404").Value = "Primera"
$ws.Range("J404").Value = 240
$ws.Range("K404").Value = 6500
$ws.Range("L404").Value = 7000
$ws.Range("M404").Value = 6750
$ws.Range("N404").Value = "`$/caja 15 unidades"
$ws.Range("O404").Value = "Región del Maule"
$ws.Range("P404").Value = 450
$ws.Range("Q404").Value = 15
$ws.Range("R404").Value = "Hortaliza"
